$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 6-cylinder group (rows 6 & 7, previously merged A6:A7) is being split
# into its two underlying (cylinder, engine-shape) groups, so the table now
# has 4 data rows instead of 3. Unmerge the old A6:A7 label cell first.
$ws.Range("A6:A7").UnMerge()

# Row 8 already carries the exact per-column formatting (styles) that every
# data row should use, so stamp it onto the other three data rows before
# writing the new values.
$ws.Range("A8:H8").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)
$ws.Range("A6:H6").PasteSpecial(-4122)
$ws.Range("A7:H7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 5: 6-cylinder, vs = 0
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 110
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 2.7475
$ws.Range("G5").Value = 0.1803122292025695

# Row 6: 4-cylinder, vs = 1
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 83.33333333333333
$ws.Range("E6").Value = 18.50225211517056
$ws.Range("F6").Value = 2.886666666666667
$ws.Range("G6").Value = 0.4911551010967242

# Row 7: 6-cylinder, vs = 1
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 112.6666666666667
$ws.Range("E7").Value = 9.291573243177568
$ws.Range("F7").Value = 3.371666666666667
$ws.Range("G7").Value = 0.1360453355809502

# Row 8 (8-cylinder) is unchanged.
